# Add files via upload
# This script recreates the changes made to the "Escenarios PE" sheet:
#  - fills in evidence video hyperlinks (column I) for rows 10-15
#  - fills in "Errores encontrados" values (column J) for rows 10-15
#  - fills in the new "Nuevo post" test scenario rows (14-15), columns E-H
#  - leaves the final selection on E16, as in the saved workbook

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Escenarios PE")

$missing = [System.Reflection.Missing]::Value

$u10 = "https://uniandes-my.sharepoint.com/:v:/g/personal/r_ramosg_uniandes_edu_co/EbDR65m1kOtOjif50MnnYhcBCr2-9jwzdCK8Bd3g17UGcQ?nav=eyJyZWZlcnJhbEluZm8iOnsicmVmZXJyYWxBcHAiOiJPbmVEcml2ZUZvckJ1c2luZXNzIiwicmVmZXJyYWxBcHBQbGF0Zm9ybSI6IldlYiIsInJlZmVycmFsTW9kZSI6InZpZXciLCJyZWZlcnJhbFZpZXciOiJNeUZpbGVzTGlua0NvcHkifX0&e=xQX5Ak"
$u11 = "https://uniandes-my.sharepoint.com/:v:/r/personal/r_ramosg_uniandes_edu_co/Documents/MISW4103-PAS/evidencias/PE02-iniciar-sesion.mp4?csf=1&web=1&nav=eyJyZWZlcnJhbEluZm8iOnsicmVmZXJyYWxBcHAiOiJPbmVEcml2ZUZvckJ1c2luZXNzIiwicmVmZXJyYWxBcHBQbGF0Zm9ybSI6IldlYiIsInJlZmVycmFsTW9kZSI6InZpZXciLCJyZWZlcnJhbFZpZXciOiJNeUZpbGVzTGlua0NvcHkifX0&e=Y5mXLI"
$u12 = "https://uniandes-my.sharepoint.com/:v:/r/personal/r_ramosg_uniandes_edu_co/Documents/MISW4103-PAS/evidencias/PE03-editar%20titulo.mp4?csf=1&web=1&nav=eyJyZWZlcnJhbEluZm8iOnsicmVmZXJyYWxBcHAiOiJPbmVEcml2ZUZvckJ1c2luZXNzIiwicmVmZXJyYWxBcHBQbGF0Zm9ybSI6IldlYiIsInJlZmVycmFsTW9kZSI6InZpZXciLCJyZWZlcnJhbFZpZXciOiJNeUZpbGVzTGlua0NvcHkifX0&e=pgEgda"
$u13 = "https://uniandes-my.sharepoint.com/:v:/r/personal/r_ramosg_uniandes_edu_co/Documents/MISW4103-PAS/evidencias/PE04-editar-titulo.mp4?csf=1&web=1&nav=eyJyZWZlcnJhbEluZm8iOnsicmVmZXJyYWxBcHAiOiJPbmVEcml2ZUZvckJ1c2luZXNzIiwicmVmZXJyYWxBcHBQbGF0Zm9ybSI6IldlYiIsInJlZmVycmFsTW9kZSI6InZpZXciLCJyZWZlcnJhbFZpZXciOiJNeUZpbGVzTGlua0NvcHkifX0&e=8gbpVr"
$u14 = "https://uniandes-my.sharepoint.com/:v:/r/personal/r_ramosg_uniandes_edu_co/Documents/MISW4103-PAS/evidencias/PE05-nuevo%20post.mp4?csf=1&web=1&nav=eyJyZWZlcnJhbEluZm8iOnsicmVmZXJyYWxBcHAiOiJPbmVEcml2ZUZvckJ1c2luZXNzIiwicmVmZXJyYWxBcHBQbGF0Zm9ybSI6IldlYiIsInJlZmVycmFsTW9kZSI6InZpZXciLCJyZWZlcnJhbFZpZXciOiJNeUZpbGVzTGlua0NvcHkifX0&e=MVwAbE"
$u15 = "https://uniandes-my.sharepoint.com/:v:/r/personal/r_ramosg_uniandes_edu_co/Documents/MISW4103-PAS/evidencias/PE06-nuevo%20post%20con%20titulo%20largo.mp4?csf=1&web=1&nav=eyJyZWZlcnJhbEluZm8iOnsicmVmZXJyYWxBcHAiOiJPbmVEcml2ZUZvckJ1c2luZXNzIiwicmVmZXJyYWxBcHBQbGF0Zm9ybSI6IldlYiIsInJlZmVycmFsTW9kZSI6InZpZXciLCJyZWZlcnJhbFZpZXciOiJNeUZpbGVzTGlua0NvcHkifX0&e=NfUvad"

# --- Column I rows 10-12: evidence video hyperlinks ------------------------
# Row 10 ends up with no cell border once the built-in "Hyperlink" style is
# applied, so clear it up front; rows 11-15 keep their original border.
$i10 = $ws.Range("I10")
$i10.Borders.LineStyle = -4142
$ws.Hyperlinks.Add($i10, $u10, $missing, $missing, $u10) | Out-Null
$ws.Hyperlinks.Add($ws.Range("I11"), $u11, $missing, $missing, $u11) | Out-Null
$ws.Hyperlinks.Add($ws.Range("I12"), $u12, $missing, $missing, $u12) | Out-Null

# --- Column J rows 10-12: "Errores encontrados" -----------------------------
$ws.Range("J10").Value = "-"
$ws.Range("J11").Value = "-"
$ws.Range("J12").Value = "-"

# --- Column I rows 13-15: evidence video hyperlinks -------------------------
$ws.Hyperlinks.Add($ws.Range("I13"), $u13, $missing, $missing, $u13) | Out-Null
$ws.Hyperlinks.Add($ws.Range("I14"), $u14, $missing, $missing, $u14) | Out-Null
$ws.Hyperlinks.Add($ws.Range("I15"), $u15, $missing, $missing, $u15) | Out-Null

# --- Column J rows 13-15: "Errores encontrados" -----------------------------
$ws.Range("J15").Value = "RI02"
$ws.Range("J13").Value = "RI01"
$ws.Range("J14").Value = "-"

# --- Row 14/15: new scenario data (Nuevo post) ------------------------------
$ws.Range("E14").Value = "Nuevo post"
$ws.Range("F14").Value = "Funcional"
$ws.Range("G14").Value = "Positivo"
$ws.Range("H14").Value = "Crear post normalmente"

$ws.Range("E15").Value = "Nuevo post"
$ws.Range("F15").Value = "Funcional"
$ws.Range("G15").Value = "Negativo"
$ws.Range("H15").Value = "Crear post usando titulo largo"

# --- Restore the selection that was active when the file was saved --------
$ws.Range("E16").Select() | Out-Null
